# Helper: force a run split at an absolute document character position by
# adding then immediately deleting a temporary bookmark there. Word (and
# this host) always separates runs cleanly at a bookmark boundary, and
# deleting the bookmark right after leaves the split in place without any
# bookmark residue.
function SplitAt($d, $pos) {
    $tmpRng = $d.Range($pos, $pos)
    $d.Bookmarks.Add("zzTempSplit", $tmpRng)
    $d.Bookmarks.Item("zzTempSplit").Delete()
}

$d = $word.ActiveDocument
$dash = [char]0x2013

# ---------------------------------------------------------------------
# Edit 1: insert the new "AGE1" coded-parameter paragraphs right before
# the "Adequacy (ZADEQ):" paragraph.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Adequacy (ZADEQ):", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$adequacyStart = $rng1.Start
$insPoint = $d.Range($adequacyStart, $adequacyStart)

$block = "AGE1`r"
$block = $block + "-9, 0 " + $dash + " 93`r"
$block = $block + "Group into buckets (13 - 33, 34 - 53, 54 - 73, 74 " + $dash + " 93)`r"
$block = $block + "Age1 (13-33)`r"
$block = $block + "Age2 (34 " + $dash + " 53)`r"
$block = $block + " Age3 (54 " + $dash + " 73)`r"
$block = $block + "Age4 (74 " + $dash + " 93)`r"
$insPoint.InsertBefore($block)

# Locate the newly inserted paragraphs by their text and apply list
# formatting (matching the existing ZADEQ-style bullet list: numId 2).
$pAge1Bucket = $d.Paragraphs.Item(1)
$rngFind = $d.Content
$rngFind.Find.Execute("-9, 0 " + $dash + " 93", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pMinus9 = $rngFind.Paragraphs.Item(1)
$pMinus9.Style = "List Paragraph"
$pMinus9.Range.ListFormat.ListLevelNumber = 1

$rngFind2 = $d.Content
$rngFind2.Find.Execute("Group into buckets", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pBuckets = $rngFind2.Paragraphs.Item(1)
$pBuckets.Style = "List Paragraph"
$pBuckets.Range.ListFormat.ListLevelNumber = 1

$rngFind3 = $d.Content
$rngFind3.Find.Execute("Age1 (13-33)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pAge1 = $rngFind3.Paragraphs.Item(1)
$pAge1.Style = "List Paragraph"
$pAge1.Range.ListFormat.ListLevelNumber = 2

$rngFind4 = $d.Content
$rngFind4.Find.Execute("Age2 (34 " + $dash + " 53)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pAge2 = $rngFind4.Paragraphs.Item(1)
$pAge2.Style = "List Paragraph"
$pAge2.Range.ListFormat.ListLevelNumber = 2

$rngFind5 = $d.Content
$rngFind5.Find.Execute("Age3 (54 " + $dash + " 73)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pAge3 = $rngFind5.Paragraphs.Item(1)
$pAge3.Style = "List Paragraph"
$pAge3.Range.ListFormat.ListLevelNumber = 2

$rngFind6 = $d.Content
$rngFind6.Find.Execute("Age4 (74 " + $dash + " 93)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pAge4 = $rngFind6.Paragraphs.Item(1)
$pAge4.Style = "List Paragraph"
$pAge4.Range.ListFormat.ListLevelNumber = 2

# Split the "-9, 0 – 93" paragraph into its 4 runs: "-9, 0 " | "–" | " " | "93"
$rngP1 = $d.Content
$rngP1.Find.Execute("-9, 0 " + $dash + " 93", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base1 = $rngP1.Start
$o1a = $base1 + 6
$o1b = $base1 + 7
$o1c = $base1 + 8
SplitAt $d $o1a
SplitAt $d $o1b
SplitAt $d $o1c

# Split the "Group into buckets (...)" paragraph into its 5 runs:
# "Group into buckets (" | "13 - 33, 34 - 53, 54 - 73, 74 " | "–" | " 93" | ")"
$rngP2 = $d.Content
$rngP2.Find.Execute("Group into buckets (13 - 33, 34 - 53, 54 - 73, 74 " + $dash + " 93)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base2 = $rngP2.Start
$o2a = $base2 + 20
$o2b = $base2 + 50
$o2c = $base2 + 51
$o2d = $base2 + 54
SplitAt $d $o2a
SplitAt $d $o2b
SplitAt $d $o2c
SplitAt $d $o2d

# ---------------------------------------------------------------------
# Edit 2: expand the "Location" note with the extra sentence about
# census regions, and move the "_GoBack" bookmark to sit right after it
# (adding a bookmark named "_GoBack" elsewhere automatically relocates
# the single allowed instance, which also satisfies edit 3 below).
# ---------------------------------------------------------------------
$rngSplit = $d.Content
$rngSplit.Find.Execute(", but with 4 categories", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPt = $rngSplit.Start
$insPoint2 = $d.Range($splitPt, $splitPt)
$insPoint2.InsertBefore(" which is for the census regions")
SplitAt $d $splitPt

$rngBm = $d.Content
$rngBm.Find.Execute(", but with 4 categories", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPt = $rngBm.Start
$bmRng = $d.Range($bmPt, $bmPt)

# ---------------------------------------------------------------------
# Edit 3: this Add call both places the bookmark in its new location and
# removes it from its old location after "Pre 1940 will be 0" (Word only
# keeps a single "_GoBack" bookmark at a time).
# ---------------------------------------------------------------------
$d.Bookmarks.Add("_GoBack", $bmRng)
